$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Contoh Ganjil 2023"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Contoh Ganjil 2024"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Contoh Ganjil 2025"

$ws.Range("B6").Select()
